$d = $word.ActiveDocument

$map = @{
    "85×23=1955" = "50×71=3550"
    "39×98=3822" = "99×71=7029"
    "33×84=2772" = "97×88=8536"
    "94×92=8648" = "67×84=5628"
    "75×65=4875" = "53×39=2067"
    "92×65=5980" = "42×94=3948"
    "12×14=168"  = "31×84=2604"
    "54×34=1836" = "45×43=1935"
    "56×52=2912" = "86×36=3096"
    "20×14=280"  = "28×57=1596"
    "88×14=1232" = "76×86=6536"
    "41×70=2870" = "93×28=2604"
    "89×42=3738" = "33×71=2343"
    "22×55=1210" = "28×36=1008"
    "78×72=5616" = "25×60=1500"
    "84×62=5208" = "32×18=576"
    "13×89=1157" = "68×48=3264"
    "84×57=4788" = "91×53=4823"
    "20×35=700"  = "17×16=272"
    "43×33=1419" = "59×26=1534"
    "91×65=5915" = "69×75=5175"
    "16×40=640"  = "58×73=4234"
    "56×26=1456" = "58×80=4640"
    "28×49=1372" = "71×33=2343"
    "89×57=5073" = "76×55=4180"
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
